# close #206: Adds support for zero-sum influencing factors
#
# Row 5 (id 1100049): the W/X influencing-factor pair becomes a hard 0/0
# (no longer 0.442/0.558).
# Row 6 (id 1100056): the W/X influencing-factor pair becomes the "DI"
# marker (shared text), matching the existing "DI" used elsewhere (e.g. F3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0

$ws.Range("W6").Value = "DI"
$ws.Range("X6").Value = "DI"

# Restore the author's last on-screen selection / scroll position.
$ws.Range("U13").Select()
